$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 1
$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(3, 1).Value = 6
$ws.Cells.Item(3, 5).Value = 4
$ws.Cells.Item(4, 1).Value = 7
$ws.Cells.Item(4, 5).Value = 5
$ws.Cells.Item(5, 1).Value = 9
$ws.Cells.Item(5, 5).Value = 7
$ws.Cells.Item(6, 1).Value = 11
$ws.Cells.Item(6, 5).Value = 8
$ws.Cells.Item(7, 1).Value = 14
$ws.Cells.Item(7, 5).Value = 10
$ws.Cells.Item(8, 1).Value = 18
$ws.Cells.Item(8, 5).Value = 13
$ws.Cells.Item(9, 1).Value = 21
$ws.Cells.Item(9, 5).Value = 15
$ws.Cells.Item(10, 1).Value = 24
$ws.Cells.Item(10, 5).Value = 17
$ws.Cells.Item(11, 1).Value = 28
$ws.Cells.Item(11, 5).Value = 19
$ws.Cells.Item(12, 1).Value = 31
$ws.Cells.Item(12, 5).Value = 21
$ws.Cells.Item(13, 1).Value = 33
$ws.Cells.Item(13, 5).Value = 22
$ws.Cells.Item(14, 1).Value = 36
$ws.Cells.Item(14, 5).Value = 25
$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 5).Value = 14
$ws.Cells.Item(16, 1).Value = 16
$ws.Cells.Item(16, 5).Value = 12
$ws.Cells.Item(17, 1).Value = 4
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(18, 1).Value = 37
$ws.Cells.Item(18, 5).Value = 26
$ws.Cells.Item(19, 1).Value = 12
$ws.Cells.Item(19, 5).Value = 9
$ws.Cells.Item(20, 1).Value = 34
$ws.Cells.Item(20, 5).Value = 23
$ws.Cells.Item(21, 1).Value = 27
$ws.Cells.Item(21, 5).Value = 20
$ws.Cells.Item(22, 1).Value = 10
$ws.Cells.Item(22, 5).Value = 11
$ws.Cells.Item(26, 1).Value = 2
$ws.Cells.Item(26, 5).Value = 2
$ws.Cells.Item(27, 1).Value = 8
$ws.Cells.Item(27, 5).Value = 6

# Rows 23-25: full row content rotates (row23<-old25, row24<-old23, row25<-old24),
# with columns A and E overridden to their new target values.
# Row 23 <- old row 25 content
$ws.Cells.Item(23, 1).Value = 17
$ws.Cells.Item(23, 2).Value = "2023-08-06"
$ws.Cells.Item(23, 3).Value = "18:30"
$ws.Cells.Item(23, 4).Value = "Série A"
$ws.Cells.Item(23, 5).Value = 18
$ws.Cells.Item(23, 6).Value = "Sun"
$ws.Cells.Item(23, 7).Value = "Home"
$ws.Cells.Item(23, 8).Value = "D"
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = "Botafogo (RJ)"
$ws.Cells.Item(23, 12).Value = 0.9
$ws.Cells.Item(23, 13).Value = 0.1
$ws.Cells.Item(23, 14).Value = 67
$ws.Cells.Item(23, 15).Value = 44759
$ws.Cells.Item(23, 16).Value = 13
$ws.Cells.Item(23, 17).Value = 4
$ws.Cells.Item(23, 18).Value = 30.8
$ws.Cells.Item(23, 19).Value = 0
$ws.Cells.Item(23, 20).Value = 0
$ws.Cells.Item(23, 21).Value = 2
$ws.Cells.Item(23, 22).Value = 0
$ws.Cells.Item(23, 23).Value = 0
$ws.Cells.Item(23, 24).Value = 0.9
$ws.Cells.Item(23, 25).Value = 0.07000000000000001
$ws.Cells.Item(23, 26).Value = -0.9
$ws.Cells.Item(23, 27).Value = -0.9
$ws.Cells.Item(23, 28).Value = 1
$ws.Cells.Item(23, 29).Value = 1
$ws.Cells.Item(23, 30).Value = 100
$ws.Cells.Item(23, 31).Value = 1
$ws.Cells.Item(23, 32).Value = 0
$ws.Cells.Item(23, 33).Value = 0
$ws.Cells.Item(23, 34).Value = 9391
$ws.Cells.Item(23, 35).Value = 2790
$ws.Cells.Item(23, 36).Value = 0
$ws.Cells.Item(23, 37).Value = 0.7
$ws.Cells.Item(23, 38).Value = 1.7
$ws.Cells.Item(23, 39).Value = 8
$ws.Cells.Item(23, 40).Value = 52
$ws.Cells.Item(23, 41).Value = 10
$ws.Cells.Item(23, 42).Value = 2
$ws.Cells.Item(23, 43).Value = 54
$ws.Cells.Item(23, 44).Value = 3
$ws.Cells.Item(23, 45).Value = 3
$ws.Cells.Item(23, 46).Value = 35
$ws.Cells.Item(23, 47).Value = 16
$ws.Cells.Item(23, 48).Value = 24
$ws.Cells.Item(23, 49).Value = 0
$ws.Cells.Item(23, 50).Value = 11
$ws.Cells.Item(23, 51).Value = 7
$ws.Cells.Item(23, 52).Value = 5
$ws.Cells.Item(23, 53).Value = 3
$ws.Cells.Item(23, 54).Value = 7
$ws.Cells.Item(23, 55).Value = 0
$ws.Cells.Item(23, 56).Value = "Cruzeiro"

# Row 24 <- old row 23 content
$ws.Cells.Item(24, 1).Value = 35
$ws.Cells.Item(24, 2).Value = "2023-09-22"
$ws.Cells.Item(24, 3).Value = "20:00"
$ws.Cells.Item(24, 4).Value = "Série A"
$ws.Cells.Item(24, 5).Value = 24
$ws.Cells.Item(24, 6).Value = "Fri"
$ws.Cells.Item(24, 7).Value = "Home"
$ws.Cells.Item(24, 8).Value = "W"
$ws.Cells.Item(24, 9).Value = 1
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = "Botafogo (RJ)"
$ws.Cells.Item(24, 12).Value = 1.1
$ws.Cells.Item(24, 13).Value = 0.5
$ws.Cells.Item(24, 14).Value = 66
$ws.Cells.Item(24, 15).Value = ""
$ws.Cells.Item(24, 16).Value = 19
$ws.Cells.Item(24, 17).Value = 5
$ws.Cells.Item(24, 18).Value = 26.3
$ws.Cells.Item(24, 19).Value = 0.05
$ws.Cells.Item(24, 20).Value = 0.2
$ws.Cells.Item(24, 21).Value = 1
$ws.Cells.Item(24, 22).Value = 0
$ws.Cells.Item(24, 23).Value = 0
$ws.Cells.Item(24, 24).Value = 1.1
$ws.Cells.Item(24, 25).Value = 0.06
$ws.Cells.Item(24, 26).Value = -0.1
$ws.Cells.Item(24, 27).Value = -0.1
$ws.Cells.Item(24, 28).Value = 2
$ws.Cells.Item(24, 29).Value = 2
$ws.Cells.Item(24, 30).Value = 100
$ws.Cells.Item(24, 31).Value = 1
$ws.Cells.Item(24, 32).Value = 0.1
$ws.Cells.Item(24, 33).Value = 0.1
$ws.Cells.Item(24, 34).Value = 11615
$ws.Cells.Item(24, 35).Value = 3657
$ws.Cells.Item(24, 36).Value = 0
$ws.Cells.Item(24, 37).Value = 0.8
$ws.Cells.Item(24, 38).Value = 0.8
$ws.Cells.Item(24, 39).Value = 17
$ws.Cells.Item(24, 40).Value = 68
$ws.Cells.Item(24, 41).Value = 6
$ws.Cells.Item(24, 42).Value = 2
$ws.Cells.Item(24, 43).Value = 62
$ws.Cells.Item(24, 44).Value = 1
$ws.Cells.Item(24, 45).Value = 6
$ws.Cells.Item(24, 46).Value = 24
$ws.Cells.Item(24, 47).Value = 11
$ws.Cells.Item(24, 48).Value = 36
$ws.Cells.Item(24, 49).Value = 1
$ws.Cells.Item(24, 50).Value = 12
$ws.Cells.Item(24, 51).Value = 8
$ws.Cells.Item(24, 52).Value = 8
$ws.Cells.Item(24, 53).Value = 5
$ws.Cells.Item(24, 54).Value = 7
$ws.Cells.Item(24, 55).Value = 0
$ws.Cells.Item(24, 56).Value = "Corinthians"

# Row 25 <- old row 24 content
$ws.Cells.Item(25, 1).Value = 21
$ws.Cells.Item(25, 2).Value = "2023-07-23"
$ws.Cells.Item(25, 3).Value = "16:00"
$ws.Cells.Item(25, 4).Value = "Série A"
$ws.Cells.Item(25, 5).Value = 16
$ws.Cells.Item(25, 6).Value = "Sun"
$ws.Cells.Item(25, 7).Value = "Home"
$ws.Cells.Item(25, 8).Value = "D"
$ws.Cells.Item(25, 9).Value = 2
$ws.Cells.Item(25, 10).Value = 2
$ws.Cells.Item(25, 11).Value = "Botafogo (RJ)"
$ws.Cells.Item(25, 12).Value = 0.9
$ws.Cells.Item(25, 13).Value = 1
$ws.Cells.Item(25, 14).Value = 37
$ws.Cells.Item(25, 15).Value = ""
$ws.Cells.Item(25, 16).Value = 7
$ws.Cells.Item(25, 17).Value = 2
$ws.Cells.Item(25, 18).Value = 28.6
$ws.Cells.Item(25, 19).Value = 0.29
$ws.Cells.Item(25, 20).Value = 1
$ws.Cells.Item(25, 21).Value = 0
$ws.Cells.Item(25, 22).Value = 0
$ws.Cells.Item(25, 23).Value = 0
$ws.Cells.Item(25, 24).Value = 0.9
$ws.Cells.Item(25, 25).Value = 0.13
$ws.Cells.Item(25, 26).Value = 1.1
$ws.Cells.Item(25, 27).Value = 1.1
$ws.Cells.Item(25, 28).Value = 6
$ws.Cells.Item(25, 29).Value = 4
$ws.Cells.Item(25, 30).Value = 50
$ws.Cells.Item(25, 31).Value = 0
$ws.Cells.Item(25, 32).Value = 1.5
$ws.Cells.Item(25, 33).Value = -0.6
$ws.Cells.Item(25, 34).Value = 4339
$ws.Cells.Item(25, 35).Value = 1664
$ws.Cells.Item(25, 36).Value = 2
$ws.Cells.Item(25, 37).Value = 0.9
$ws.Cells.Item(25, 38).Value = 0.6
$ws.Cells.Item(25, 39).Value = 7
$ws.Cells.Item(25, 40).Value = 16
$ws.Cells.Item(25, 41).Value = 3
$ws.Cells.Item(25, 42).Value = 1
$ws.Cells.Item(25, 43).Value = 23
$ws.Cells.Item(25, 44).Value = 1
$ws.Cells.Item(25, 45).Value = 5
$ws.Cells.Item(25, 46).Value = 9
$ws.Cells.Item(25, 47).Value = 8
$ws.Cells.Item(25, 48).Value = 12
$ws.Cells.Item(25, 49).Value = 3
$ws.Cells.Item(25, 50).Value = 13
$ws.Cells.Item(25, 51).Value = 10
$ws.Cells.Item(25, 52).Value = 7
$ws.Cells.Item(25, 53).Value = 2
$ws.Cells.Item(25, 54).Value = 18
$ws.Cells.Item(25, 55).Value = 0
$ws.Cells.Item(25, 56).Value = "Santos"
